# Update segmentation with new timing (LWP2_0011 lab timing)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("In Lab")

# --- Update Start Time values (column B) with new, second-precision timings ---
$ws.Range("B8").Value  = 0.4962847222222222
$ws.Range("B9").Value  = 0.49778935185185186
$ws.Range("B10").Value = 0.49914351851851851
$ws.Range("B11").Value = 0.50788194444444446
$ws.Range("B12").Value = 0.5091782407407407
$ws.Range("B13").Value = 0.51009259259259265
$ws.Range("B14").Value = 0.51765046296296291
$ws.Range("B15").Value = 0.5191782407407407
$ws.Range("B16").Value = 0.51962962962962966
$ws.Range("B17").Value = 0.52219907407407407
$ws.Range("B18").Value = 0.53001157407407407
$ws.Range("B19").Value = 0.53105324074074078
$ws.Range("B20").Value = 0.53140046296296295
$ws.Range("B21").Value = 0.53619212962962959
$ws.Range("B22").Value = 0.53746527777777775
$ws.Range("B23").Value = 0.53857638888888892
$ws.Range("B24").Value = 0.53918981481481476

# --- Re-format the Start Time column (B7:B24) to show seconds (h:mm -> h:mm:ss) ---
$ws.Range("B7:B24").NumberFormat = "h:mm:ss"

# --- Move the active selection from K16 to B25, and drop the frozen/scrolled
#     top-left cell override (back to default A1) ---
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("B25").Select()
